$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-26 02:16:40"
$wsOverview.Range("G3").Value = "2016-08-26 02:16:40"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-26 02:16:35"
$wsZhCn.Range("K2").Value = "2016-08-26 02:16:51"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-26 02:16:35"
$wsZhCn.Range("K3").Value = "2016-08-26 02:16:51"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-26 02:16:40"
$wsDeDe.Range("K2").Value = "2016-08-26 02:16:58"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-26 02:16:40"
$wsDeDe.Range("K3").Value = "2016-08-26 02:16:58"
